# "Fixed link in system overview"
#
# Slide 10 ("Electronics Tooling"), shape id=3 ("Content Placeholder 2"):
#   1. The text box's autofit was set to "Shrink text on overflow" with a
#      10% line-spacing reduction baked in (<a:normAutofit lnSpcReduction="10000"/>).
#      The fix re-lets PowerPoint recompute the fit, dropping the explicit
#      reduction (<a:normAutofit/>).
#   2. A stray duplicate of the "Latest EPS" hyperlink (+ trailing run) that
#      had been pasted into the last bullet is removed, leaving an empty
#      paragraph behind.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(10)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame

# --- 1. Let PowerPoint re-fit the text instead of forcing a fixed shrink ---
$tf.AutoSize = 2   # ppAutoSizeTextToFitShape -> <a:normAutofit/> (no lnSpcReduction)

# --- 2. Remove the duplicated hyperlink run (+ trailing space run) from the
#        last "sharepoint:" bullet, leaving the paragraph mark behind ---
$tr  = $tf.TextRange
$url = "https://sotonac.sharepoint.com/:f:/r/teams/UniversityofSouthamptonSmallSatelliteUoS3/FEE%20GDP%20202021/EPS/Power_Boardv3.2_270321?csf=1&web=1&e=u4JDre"

$full = $tr.Text
$firstIdx = $full.IndexOf($url)
$secondIdx = $full.IndexOf($url, $firstIdx + 1)

if ($secondIdx -ge 0) {
    $crIdx = $full.IndexOf([char]13, $secondIdx)
    if ($crIdx -lt 0) { $crIdx = $full.Length }
    $len = $crIdx - $secondIdx
    $start = $secondIdx + 1   # TextRange.Characters is 1-indexed
    $tr.Characters($start, $len).Text = ""
}
